# Updates the cryptos list sheet with refreshed price / volume data,
# matching the commit "Updated cryptos list ... with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as plain text, even when it looks like a number
# (e.g. "326.98"), without leaving a lingering custom cell style behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Simple Price (D) / Volume(1h) (E) updates, keyed by row number ---
$updates = @{
    2  = @{ D = "28.736.35"; E = "  +1.65%  " }
    3  = @{ D = "1.869.57";  E = "  +1.55%  " }
    4  = @{ E = "  +0.12%  " }
    5  = @{ D = "326.98";    E = "  -1.67%  " }
    6  = @{ D = "1.006";     E = "  +0.30%  " }
    8  = @{ D = "0.3923";    E = "  +1.07%  " }
    9  = @{ D = "0.07919";   E = "  +0.63%  " }
    10 = @{ D = "0.9739";    E = "  +0.59%  " }
    11 = @{ D = "22.36";     E = "  +1.76%  " }
    12 = @{ D = "1.878.90";  E = "  +0.63%  " }
    13 = @{ D = "5.736";     E = "  -0.61%  " }
    14 = @{ D = "6.956";     E = "  +0.30%  " }
    15 = @{ D = "0.06925";   E = "  +0.53%  " }
    16 = @{ D = "88.80";     E = "  +1.79%  " }
    17 = @{ D = "1.007";     E = "  +0.24%  " }
    18 = @{ E = "  +0.90%  " }
    19 = @{ E = "  +0.15%  " }
    20 = @{ D = "1.004";     E = "  -0.02%  " }
    21 = @{ D = "28.776.24"; E = "  +1.70%  " }
    22 = @{ D = "5.336";     E = "  -0.40%  " }
    23 = @{ D = "11.10";     E = "  -0.08%  " }
    24 = @{ D = "2.134";     E = "  -1.47%  " }
    25 = @{ D = "2.181.41";  E = "  +4.27%  " }
    26 = @{ D = "155.44";    E = "  +1.25%  " }
    27 = @{ D = "19.30";     E = "  -0.05%  " }
    28 = @{ D = "5.776";     E = "  -2.09%  " }
    29 = @{ D = "1.997";     E = "  +0.81%  " }
    30 = @{ D = "119.31";    E = "  +1.96%  " }
    31 = @{ D = "0.09359";   E = "  +0.18%  " }
    32 = @{ D = "0.9410";    E = "  -0.84%  " }
    33 = @{ D = "5.335";     E = "  +0.32%  " }
    34 = @{ E = "  +1.22%  " }
    35 = @{ D = "3.348";     E = "  -3.14%  " }
    36 = @{ D = "0.05840";   E = "  -3.77%  " }
    37 = @{ D = "0.02117";   E = "  -2.79%  " }
    38 = @{ D = "1.157";     E = "  +0.07%  " }
    39 = @{ D = "7.899";     E = "  +3.96%  " }
    40 = @{ D = "0.5666";    E = "  +0.36%  " }
    41 = @{ D = "9.984";     E = "  -0.62%  " }
    42 = @{ D = "0.1780";    E = "  -0.36%  " }
    43 = @{ D = "0.07334";   E = "  +4.28%  " }
    44 = @{ D = "2.255";     E = "  -5.55%  " }
    47 = @{ D = "1.143";     E = "  -6.39%  " }
    48 = @{ D = "1.854";     E = "  +0.20%  " }
    49 = @{ D = "113.99";    E = "  +0.55%  " }
    50 = @{ D = "2.358";     E = "  +1.62%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        Set-TextValue $ws.Range("D$row") $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$row").Value = $vals["E"]
    }
}

# --- Rows 45 & 46 swap places (Decentraland <-> EnergySwap) with new data ---
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D45") "11.71"
$ws.Range("E45").Value = "  -0.12%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D46") "0.5338"
$ws.Range("E46").Value = "  +0.17%  "

# --- Row 51: EOS replaced by PaxDollar ---
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D51") "1.006"
$ws.Range("E51").Value = "  +0.28%  "
